# Ignore sheets whose name starts with "!" and ignore empty sheets.
# This edit:
#   - Renames "Sheet3" to "!Sheet3" so downstream tooling skips it.
#   - Fills "!Sheet3" with some data so it is no longer empty (and
#     makes it the active sheet/tab, matching the author's workbook view).
#   - Clears the now-stale "active tab" selection marker from "First".

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("First")
$sheet3 = $wb.Worksheets.Item("Sheet3")

# Rename the third sheet so it gets ignored by the "!" convention.
$sheet3.Name = "!Sheet3"

# Populate the previously-empty sheet with a small 2x2 table.
$sheet3.Range("A1").Value = "Tach"
$sheet3.Range("B1").Value = "Moin"
$sheet3.Range("A2").Value = "d"
$sheet3.Range("B2").Value = "d"

# Restore the first sheet's own selection (no longer the active tab).
$sheet1.Range("A2").Select() | Out-Null

# "!Sheet3" becomes the active tab, with its own selection.
$sheet3.Activate()
$sheet3.Range("C13").Select() | Out-Null
